$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 25.47000000000054
$ws.Range("H2").Value = [double]"2.453531546132943e-16"
$ws.Range("K2").Value = 51.19227769964967
$ws.Range("L2").Value = "[41.38143799806999, 61.00311740122934]"
$ws.Range("O2").Value = 1.50318447288881
$ws.Range("P2").Value = "[1.2893423303021168, 1.7170266154755032]"
$ws.Range("S2").Value = 67.50247185585825
$ws.Range("T2").Value = "[61.28657165031038, 73.71837206140611]"
$ws.Range("W2").Value = 19.37657657657699
$ws.Range("X2").Value = 18.50972972973012
$ws.Range("Y2").Value = 20.24342342342386

# Row 3
$ws.Range("E3").Value = 25.39000000000053
$ws.Range("H3").Value = [double]"2.453531546132943e-16"
$ws.Range("K3").Value = 47.769988475892
$ws.Range("L3").Value = "[36.55825030822983, 58.98172664355417]"
$ws.Range("M3").Value = [double]"1.554312234475219e-15"
$ws.Range("N3").Value = [double]"1.554312234475219e-15"
$ws.Range("O3").Value = 1.037763339023655
$ws.Range("P3").Value = "[0.786184347745194, 1.2893423303021168]"
$ws.Range("Q3").Value = [double]"9.547918011776346e-15"
$ws.Range("R3").Value = [double]"9.547918011776346e-15"
$ws.Range("S3").Value = 63.70868812473351
$ws.Range("T3").Value = "[57.63691483031053, 69.78046141915648]"
$ws.Range("W3").Value = 21.1964564564569
$ws.Range("X3").Value = 20.17983983984027
$ws.Range("Y3").Value = 22.21307307307353
